$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 and Q1, matching the style of the other header cells ---
$headerRange = $ws.Range("P1:Q1")

# Match style of an existing header cell (bold font, centered/top aligned, thin border all around)
# by copying directly from O1 into the destination range (no clipboard/selection paste involved).
$ws.Range("O1").Copy($headerRange)

# Now set the actual values for the new header cells.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Body rows 2..25 ---
# For each row: swap I/K and M/O values (1<->2), and add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2 (new)
}
